# Applies the "fixing according to feedback" edit to the Algorithm list.
$d = $word.ActiveDocument
$LDQ = [char]0x201C   # left double quote "
$RDQ = [char]0x201D   # right double quote "
$NDASH = [char]0x2013 # en dash

# --- Do all simple whole-paragraph text replacements FIRST, while every
#     search string is still unique in the document (before any of the
#     new text below introduces look-alike substrings). ---

# 1) "birth" prompt paragraph (was 3 runs)
$d.Content.Find.Execute(
    "Prompt user for the number of seconds between each birth, death, and immigration in a country.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Ask user to input seconds between birth.", 2)

# 2) "death" prompt paragraph (was 3 runs)
$d.Content.Find.Execute(
    "Prompt user to input the current population of a country and the number of years in the future they want to know the population.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Ask user to input seconds between death.", 2)

# 3) "seconds in a year" paragraph becomes the "immigration" prompt
$d.Content.Find.Execute(
    "Calculate the number of seconds in a year and store it in a variable.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Ask user to input seconds between immigration.", 2)

# 4) "Keep results..." (ilvl 1) paragraph -> "Output "The total population has increased""
$d.Content.Find.Execute(
    "Keep results as floats until after the calculations have been performed.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    ("Output " + $LDQ + "The total population has increased" + $RDQ), 2)

# 5) "Output whether..." (ilvl 0) paragraph -> "Else:"
$d.Content.Find.Execute(
    "Output whether the population will have increased or decreased.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Else:", 2)

# --- Now insert the new list paragraphs, anchored by paragraph Index. ---

# Locate "Ask user to input seconds between immigration." (was step 3 above).
$immigrationIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Ask user to input seconds between immigration.`r") {
        $immigrationIdx = $i
        break
    }
}

$p = $d.Paragraphs.Item($immigrationIdx)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($p.Index + 1)
$p.Range.Text = "Ask user to input current population."

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($p.Index + 1)
$p.Range.Text = "Ask user to input years in future."

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($p.Index + 1)
$p.Range.Text = "Calculate the number of seconds in a year and store it in a variable. "

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($p.Index + 1)
$p.Range.Text = ("Calculate the population change, and store result in a variable. " +
    "Use the formula " + $LDQ +
    "(Seconds per year/seconds between birth + seconds per year/seconds between immigration " +
    $NDASH + " seconds per year/seconds between death) * years in future" + $RDQ +
    " to calculate for the population change.")

# Append extra detail to the (still untouched) "Calculate the expected population" paragraph.
$expIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Calculate the expected population in the future and output it as an integer.`r") {
        $expIdx = $i
        break
    }
}
$expPara = $d.Paragraphs.Item($expIdx)
$endRng = $d.Range($expPara.Range.End - 1, $expPara.Range.End - 1)
$endRng.InsertAfter(" Use formula " + $LDQ + "(current population + population change)" + $RDQ +
    " to calculate future population. Keep results as floats until after the calculations have been performed.")

# Insert "If future population >= current population:" right after it.
$expPara = $d.Paragraphs.Item($expIdx)
$expPara.Range.InsertParagraphAfter()
$ifPara = $d.Paragraphs.Item($expPara.Index + 1)
$ifPara.Range.Text = "If future population >= current population:"

# Insert the final "decreased" paragraph (ilvl 1) right after "Else:".
$elseIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Else:`r") {
        $elseIdx = $i
        break
    }
}
$elsePara = $d.Paragraphs.Item($elseIdx)
$elsePara.Range.InsertParagraphAfter()
$decPara = $d.Paragraphs.Item($elsePara.Index + 1)
$decPara.Range.ListFormat.ListLevelNumber = 2
$decPara.Range.Text = "Output " + $LDQ + "The total population has decreased." + $RDQ

Write-Output "done"
